$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "consequents_length" in H1, matching the style of the other
# header cells (bold font + thin border + centered alignment), by copying the
# format from the existing G1 header cell.
$ws.Range("H1").Value = "consequents_length"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill the new column with the consequents_length value (1) for every data row.
$ws.Range("H2:H10").Value = 1
